{"js": "// Remove the leftover \"<< MAYBE DELETE THIS - IT IS REPEATIVE >>\" review\n// comment paragraph, the repetitive recap sentence right after it, and the\n// blank paragraph that followed it - three whole paragraphs in total -\n// right after the \"Machine Learning and Analysis Recap:\" heading.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(\"MAYBE DELETE THIS\") !== -1) {\n        targetIndex = i;\n        break;\n    }\n}\n\nif (targetIndex !== -1) {\n    // Delete the marker paragraph plus the next two paragraphs (the\n    // repeated recap sentence and the blank line that followed it).\n    paragraphs.items[targetIndex].delete();\n    paragraphs.items[targetIndex + 1].delete();\n    paragraphs.items[targetIndex + 2].delete();\n    await context.sync();\n}\n", "ps1": "# Remove the leftover \"<< MAYBE DELETE THIS - IT IS REPEATIVE >>\" review\n# comment paragraph, the repetitive recap sentence right after it, and the\n# blank paragraph that followed it - three whole paragraphs in total -\n# right after the \"Machine Learning and Analysis Recap:\" heading.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*MAYBE DELETE THIS*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Deleting the paragraph's Range removes the paragraph mark too, so the\n    # paragraph that used to follow \"$target\" slides into its place. Doing\n    # this three times removes the marker paragraph plus the next two\n    # paragraphs (the repeated recap sentence and the blank line after it).\n    for ($i = 0; $i -lt 3; $i++) {\n        $target.Range.Delete()\n    }\n}\n"}
